# Apply updated cryptocurrency price/volume values to the worksheet.
# Values are forced to be stored as literal text (matching the source
# inlineStr cells) by prefixing with an apostrophe and then resetting
# the cell style to "Normal" so no numeric formatting/quote-prefix style
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "30.102.13"
Set-TextValue "E2" "  -0.63%  "
Set-TextValue "D3" "1.856.25"
Set-TextValue "E3" "  -0.73%  "
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "233.57"
Set-TextValue "E5" "  -0.76%  "
Set-TextValue "E6" "  +0.02%  "
Set-TextValue "D7" "0.4689"
Set-TextValue "D8" "42.83"
Set-TextValue "E8" "  -0.25%  "
Set-TextValue "D9" "0.2825"
Set-TextValue "E9" "  -1.74%  "
Set-TextValue "D10" "0.06443"
Set-TextValue "E10" "  -2.09%  "
Set-TextValue "D11" "20.96"
Set-TextValue "E11" "  -4.02%  "
Set-TextValue "D12" "0.07737"
Set-TextValue "E12" "  -3.54%  "
Set-TextValue "D13" "1.869.22"
Set-TextValue "E13" "  -0.07%  "
Set-TextValue "D14" "93.37"
Set-TextValue "E14" "  -3.98%  "
Set-TextValue "D15" "0.6779"
Set-TextValue "E15" "  -1.18%  "
Set-TextValue "D16" "5.041"
Set-TextValue "E16" "  -1.71%  "
Set-TextValue "D17" "265.43"
Set-TextValue "E17" "  -1.57%  "
Set-TextValue "D18" "30.092.84"
Set-TextValue "E18" "  -0.60%  "
Set-TextValue "D19" "13.30"
Set-TextValue "E19" "  -5.42%  "
Set-TextValue "D20" "0.000007562"
Set-TextValue "E20" "  -1.51%  "
Set-TextValue "E21" "  +0.01%  "
Set-TextValue "D22" "2.112.60"
Set-TextValue "E22" "  -0.05%  "
Set-TextValue "E23" "  +0.02%  "
Set-TextValue "D24" "5.137"
Set-TextValue "E24" "  -2.76%  "
Set-TextValue "D25" "6.090"
Set-TextValue "E25" "  -2.08%  "
Set-TextValue "D26" "9.295"
Set-TextValue "E26" "  -1.46%  "
Set-TextValue "D27" "165.21"
Set-TextValue "E27" "  -1.82%  "
Set-TextValue "D28" "18.45"
Set-TextValue "E28" "  -2.42%  "
Set-TextValue "D29" "1.879"
Set-TextValue "E29" "  -3.56%  "
Set-TextValue "E30" "  -0.65%  "
Set-TextValue "D31" "0.09831"
Set-TextValue "E31" "  -0.54%  "
Set-TextValue "D32" "1.449"
Set-TextValue "E32" "  -1.02%  "
Set-TextValue "D33" "4.186"
Set-TextValue "E33" "  -4.46%  "
Set-TextValue "D34" "3.969"
Set-TextValue "E34" "  -2.69%  "
Set-TextValue "D35" "0.04639"
Set-TextValue "E35" "  -1.38%  "
Set-TextValue "E36" "  -2.04%  "
Set-TextValue "D37" "0.6874"
Set-TextValue "E37" "  -1.94%  "
Set-TextValue "D39" "0.01830"
Set-TextValue "E39" "  -2.42%  "
Set-TextValue "D40" "2.718"
Set-TextValue "E40" "  +3.49%  "
Set-TextValue "D41" "6.276"
Set-TextValue "E41" "  -0.36%  "
Set-TextValue "D42" "70.51"
Set-TextValue "E42" "  -2.63%  "
Set-TextValue "E43" "  +0.04%  "
Set-TextValue "D44" "0.8316"
Set-TextValue "E44" "  -1.37%  "
Set-TextValue "D45" "1.869"
Set-TextValue "E45" "  -4.35%  "
Set-TextValue "D46" "101.66"
Set-TextValue "E46" "  -1.44%  "
Set-TextValue "D47" "0.4036"
Set-TextValue "E47" "  -3.13%  "
Set-TextValue "D48" "9.133"
Set-TextValue "E48" "  -0.62%  "
Set-TextValue "D49" "6.924"
Set-TextValue "E49" "  -1.96%  "
Set-TextValue "D50" "918.91"
Set-TextValue "E50" "  -0.50%  "
Set-TextValue "D51" "34.05"
Set-TextValue "E51" "  -1.39%  "
